$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values - C2 and E2 are cleared, B2 and D2 updated
$ws.Range("B2").Value = 5.9186518994940718
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 3.7323335046658919
$ws.Range("E2").ClearContents()

# Row 3 values updated
$ws.Range("B3").Value = 4.8271078699636059
$ws.Range("C3").Value = 5.6375100864256718
$ws.Range("D3").Value = 3.5820636464601581
$ws.Range("E3").Value = 8.7406576949142938

# Update selection to match new sqref
$ws.Range("B1:E3").Select()
